# Update betting odds / correct-score cells on the single worksheet to match
# the latest FlashScore data pull (2025-02-18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.88
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1.91
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("U3").Value = 4.1
$ws.Range("V3").Value = 1.23
$ws.Range("AA3").Value = 2.2
$ws.Range("AB3").Value = 1.62
$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 10
$ws.Range("AI3").Value = 6.5
$ws.Range("AL3").Value = 81
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 1.91
$ws.Range("L4").Value = 3.6
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.88
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 11
$ws.Range("AG4").Value = 26
$ws.Range("AN4").Value = 7
$ws.Range("AO4").Value = 12
$ws.Range("AQ4").Value = 29
$ws.Range("O5").Value = 1.72
$ws.Range("P5").Value = 2
$ws.Range("W5").Value = 5.8
$ws.Range("Y5").Value = 1.75
$ws.Range("Z5").Value = 1.98
$ws.Range("AB5").Value = 1.55
$ws.Range("AL5").Value = 150
$ws.Range("AS5").Value = 60
$ws.Range("G6").Value = 2.92
$ws.Range("I6").Value = 2.75
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 1.85
$ws.Range("L6").Value = 3.5
$ws.Range("M6").Value = 1.15
$ws.Range("N6").Value = 4.7
$ws.Range("S6").Value = 2.85
$ws.Range("T6").Value = 1.37
$ws.Range("Y6").Value = 1.57
$ws.Range("Z6").Value = 2.25
$ws.Range("AA6").Value = 2.22
$ws.Range("AB6").Value = 1.6
$ws.Range("AC6").Value = 6.4
$ws.Range("AD6").Value = 13
$ws.Range("AI6").Value = 4.7
$ws.Range("AK6").Value = 19
$ws.Range("AN6").Value = 5.7
$ws.Range("AO6").Value = 11.75
$ws.Range("AQ6").Value = 35
$ws.Range("G7").Value = 2.4
$ws.Range("H7").Value = 2.88
$ws.Range("I7").Value = 3.3
$ws.Range("J7").Value = 3.1
$ws.Range("K7").Value = 1.91
$ws.Range("Y7").Value = 1.62
$ws.Range("Z7").Value = 2.2
$ws.Range("AD7").Value = 9.5
$ws.Range("AF7").Value = 23
$ws.Range("AN7").Value = 6.5
$ws.Range("AP7").Value = 13
$ws.Range("N8").Value = 9
$ws.Range("G9").Value = 2.35
$ws.Range("H9").Value = 3.5
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 2.1
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("AA9").Value = 1.8
$ws.Range("AB9").Value = 1.95
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 9.5
$ws.Range("AK9").Value = 15
$ws.Range("AN9").Value = 8.5
$ws.Range("G10").Value = 3.75
$ws.Range("H10").Value = 3.25
$ws.Range("J10").Value = 4.75
$ws.Range("K10").Value = 1.95
$ws.Range("L10").Value = 2.88
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.5
$ws.Range("P10").Value = 2.5
$ws.Range("Q10").Value = 1.93
$ws.Range("R10").Value = 1.93
$ws.Range("S10").Value = 2.5
$ws.Range("T10").Value = 1.5
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 1.17
$ws.Range("Y10").Value = 1.57
$ws.Range("Z10").Value = 2.25
$ws.Range("AA10").Value = 2.1
$ws.Range("AB10").Value = 1.67
$ws.Range("AG10").Value = 41
$ws.Range("AI10").Value = 7
$ws.Range("AL10").Value = 81
$ws.Range("AN10").Value = 5.5
$ws.Range("AR10").Value = 21
$ws.Range("AS10").Value = 41
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 4.75
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 5.5
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 1.78
$ws.Range("R11").Value = 2.1
$ws.Range("S11").Value = 2.35
$ws.Range("T11").Value = 1.57
$ws.Range("W11").Value = 4.33
$ws.Range("X11").Value = 1.2
$ws.Range("AA11").Value = 2.1
$ws.Range("AB11").Value = 1.67
$ws.Range("AC11").Value = 5.5
$ws.Range("AI11").Value = 7
$ws.Range("AM11").Value = 501
$ws.Range("AN11").Value = 10
$ws.Range("AO11").Value = 23
$ws.Range("AP11").Value = 17
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 2.75
$ws.Range("S12").Value = 2.25
$ws.Range("T12").Value = 1.62
$ws.Range("K13").Value = 2.38
$ws.Range("Y13").Value = 1.36
$ws.Range("Z13").Value = 3
$ws.Range("AA13").Value = 2
$ws.Range("AB13").Value = 1.75
$ws.Range("AD13").Value = 6.5
$ws.Range("AF13").Value = 9
$ws.Range("AH13").Value = 29
$ws.Range("AI13").Value = 11
$ws.Range("AK13").Value = 21
$ws.Range("AL13").Value = 67
$ws.Range("AM13").Value = 351
$ws.Range("AP13").Value = 23
$ws.Range("G14").Value = 2.15
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 3.2
$ws.Range("J14").Value = 2.88
$ws.Range("L14").Value = 3.6
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 8
$ws.Range("AC14").Value = 8
$ws.Range("AD14").Value = 11
$ws.Range("AE14").Value = 9.5
$ws.Range("AF14").Value = 21
$ws.Range("AG14").Value = 19
$ws.Range("AP14").Value = 12
$ws.Range("AQ14").Value = 34
$ws.Range("AR14").Value = 26
$ws.Range("AS14").Value = 34
